$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Capture the existing number formats used so we can swap them between B7/B8.
$fmtLastRow = $ws.Range("B7").NumberFormat   # currently "YYYY-MM-DD" (style index 3)
$fmtNormal  = $ws.Range("B6").NumberFormat   # "YYYY-MM-DD HH:MM:SS" (style index 2)

# B7 reverts to the standard date format used by the other data rows, since it is
# no longer the last row once row 8 is appended.
$ws.Range("B7").NumberFormat = $fmtNormal

# New row 8 data (raw/clean SSA data for June 7th, 2020)
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 43989
$ws.Range("C8").Value = 117103
$ws.Range("D8").Value = 173975
$ws.Range("E8").Value = 45317
$ws.Range("F8").Value = 13699
$ws.Range("G8").Value = 33.48

# A8 matches the style used by the other cells in column A (bold/border/centered header style).
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# B8 takes on the "last row" date style that B7 used to have.
$ws.Range("B8").NumberFormat = $fmtLastRow
